# Auto-generated edit script: apply value updates per the commit diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$wsExh = $wb.Worksheets.Item("展览")
$wsExh.Range("F2").Value = 6780
$wsExh.Range("F3").Value = 815
$wsExh.Range("F6").Value = 9
$wsExh.Range("F7").Value = 725
$wsExh.Range("F8").Value = 725
$wsExh.Range("F9").Value = 15
$wsExh.Range("F12").Value = 1107
$wsExh.Range("F13").Value = 864
$wsExh.Range("F14").Value = 14
$wsExh.Range("F16").Value = 1016
$wsExh.Range("F17").Value = 1344
$wsExh.Range("F18").Value = 46
$wsExh.Range("F19").Value = 119
$wsExh.Range("F20").Value = 535
$wsExh.Range("F21").Value = 3
$wsExh.Range("F22").Value = 560
$wsExh.Range("F23").Value = 11
$wsExh.Range("F25").Value = 366
$wsExh.Range("F26").Value = 1062
$wsExh.Range("F27").Value = 1494
$wsExh.Range("F28").Value = 725
$wsExh.Range("F29").Value = 535
$wsExh.Range("F30").Value = 461
$wsExh.Range("F31").Value = 456
$wsExh.Range("F34").Value = 1130
$wsExh.Range("F35").Value = 264
$wsExh.Range("F36").Value = 2360
$wsExh.Range("F37").Value = 263
$wsExh.Range("F38").Value = 1254
$wsExh.Range("F39").Value = 445
$wsExh.Range("F40").Value = 63
$wsExh.Range("F41").Value = 3869

# --- Sheet: 演出 ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("G2").Value = "不可售"
$wsShow.Range("G3").Value = "不可售"
$wsShow.Range("F11").Value = 159
$wsShow.Range("F12").Value = 642
$wsShow.Range("F17").Value = 394

# --- Sheet: 本地生活 ---
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F4").Value = 1249
$wsLocal.Range("F5").Value = 1634
$wsLocal.Range("F7").Value = 136
$wsLocal.Range("F8").Value = 965

# --- Sheet: 全部类型 ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 1249
$wsAll.Range("F4").Value = 1634
$wsAll.Range("F6").Value = 136
$wsAll.Range("F7").Value = 965
$wsAll.Range("F8").Value = 6780
$wsAll.Range("F10").Value = 815
$wsAll.Range("F13").Value = 9
$wsAll.Range("F14").Value = 725
$wsAll.Range("F15").Value = 725
$wsAll.Range("F18").Value = 1107
$wsAll.Range("F19").Value = 864
$wsAll.Range("F22").Value = 159
$wsAll.Range("F23").Value = 159
$wsAll.Range("F25").Value = 1016
$wsAll.Range("F26").Value = 1344
$wsAll.Range("F27").Value = 46
$wsAll.Range("F28").Value = 119
$wsAll.Range("F29").Value = 535
$wsAll.Range("F30").Value = 560
$wsAll.Range("F33").Value = 366
$wsAll.Range("F34").Value = 1062
$wsAll.Range("F35").Value = 1494
$wsAll.Range("F36").Value = 725
$wsAll.Range("F37").Value = 535
$wsAll.Range("F38").Value = 461
$wsAll.Range("F39").Value = 456
$wsAll.Range("F43").Value = 1130
$wsAll.Range("F44").Value = 264
$wsAll.Range("F45").Value = 2360
$wsAll.Range("F49").Value = 1254
$wsAll.Range("F50").Value = 445
$wsAll.Range("F51").Value = 3869
